$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that no longer hold content (moved away / removed) ---
$ws.Range("H3").Clear()
$ws.Range("O5").Clear()
$ws.Range("O6").Clear()
$ws.Range("O7").Clear()
$ws.Range("O8").Clear()
$ws.Range("O9").Clear()
$ws.Range("H10").Clear()
$ws.Range("H11").Clear()

# --- Column P: statistics / reports block (previously in columns H/O) ---
$ws.Range("P3").Value = "статистика"
$ws.Range("P3").Font.Bold = $false

$ws.Range("P4").Value = "DepositExtractor"
$ws.Range("P4").Font.Bold = $true

$ws.Range("P5").Value = "Находит все операции по данному счету"
$ws.Range("P5").Font.Bold = $false

$ws.Range("P6").Value = "и составляет таблицу ежедневных остатков"
$ws.Range("P6").Font.Bold = $false

$ws.Range("P7").Value = "и общие суммы взносов, процентов, расходов"
$ws.Range("P7").Font.Bold = $false

# --- Column H: reporter block (new content) ---
$ws.Range("H4").Value = "отчеты"
$ws.Range("H4").Font.Bold = $false

$ws.Range("H5").Value = "DepositReporter"
$ws.Range("H5").Font.Bold = $true

$ws.Range("H6").Value = "составляет List<String> для отчета"
$ws.Range("H6").Font.Bold = $false

$ws.Range("H7").Value = "DepositExcelReporter"
$ws.Range("H7").Font.Bold = $true

$ws.Range("H8").Value = "составляет файл экселя"
$ws.Range("H8").Font.Bold = $false

# --- Column B: requirements list ---
$ws.Range("B7").Value = "нужен прогноз по месяцу"
$ws.Range("B7").Font.Bold = $false

$ws.Range("B8").Value = " и до конца депозита"
$ws.Range("B8").Font.Bold = $false

# --- Columns K/M: calculation / aggregation block (new) ---
$ws.Range("K8").Value = "агрегирование"
$ws.Range("K8").Font.Bold = $false

$ws.Range("M8").Value = "расчет"
$ws.Range("M8").Font.Bold = $false

$ws.Range("K9").Value = "DepositCalculationAggregator"
$ws.Range("K9").Font.Bold = $true

$ws.Range("M9").Value = "DepositCalculator"
$ws.Range("M9").Font.Bold = $true

$ws.Range("K10").Value = "определяет какой период "
$ws.Range("K10").Font.Bold = $false

$ws.Range("M10").Value = "расчитывает проценты по вкладу"
$ws.Range("M10").Font.Bold = $false

$ws.Range("K11").Value = "уже оплачен, какой нет"
$ws.Range("K11").Font.Bold = $false

$ws.Range("M11").Value = "за каждый день"
$ws.Range("M11").Font.Bold = $false

$ws.Range("K12").Value = "суммированием определяет "
$ws.Range("K12").Font.Bold = $false

$ws.Range("K13").Value = "проценты за опред период"
$ws.Range("K13").Font.Bold = $false

# --- Column widths (best achievable match given engine rounding) ---
$ws.Range("B1").ColumnWidth = 25.307291666666668
$ws.Range("J1").ColumnWidth = 17.592447916666668
$ws.Range("K1").ColumnWidth = 27.166666666666668
$ws.Range("L1").ColumnWidth = 9.307291666666666
$ws.Range("M1").ColumnWidth = 21.736979166666668
$ws.Range("P1").ColumnWidth = 17.022135416666668

# --- Selection matches author's final cursor position ---
$ws.Range("H4").Select()
